$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-9, columns A:T
# Row 2 and Row 3 are rewritten; rows 4-9 are newly added.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl5"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.708386
$ws.Range("H2").Value = 5.125158
$ws.Range("I2").Value = 0.09367635209466295
$ws.Range("J2").Value = 0.09367635209466295
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 114.5606336666667
$ws.Range("N2").Value = 343.681901
$ws.Range("O2").Value = 0.7368570786832789
$ws.Range("P2").Value = 0.736857078683279
$ws.Range("Q2").Value = 195.713782707262
$ws.Range("R2").Value = 1761.424044365358
$ws.Range("S2").Value = 0.0690260831461796
$ws.Range("T2").Value = 0.0690260831461796

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl5"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.708386
$ws.Range("H3").Value = 5.125158
$ws.Range("I3").Value = 0.09367635209466295
$ws.Range("J3").Value = 0.09367635209466295
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 40.91135266666667
$ws.Range("N3").Value = 122.734058
$ws.Range("O3").Value = 0.263142921316721
$ws.Range("P3").Value = 0.2631429213167211
$ws.Range("Q3").Value = 69.892382136796
$ws.Range("R3").Value = 629.031439231164
$ws.Range("S3").Value = 0.02465026894848335
$ws.Range("T3").Value = 0.02465026894848335

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl5"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.842846333333334
$ws.Range("H4").Value = 11.528539
$ws.Range("I4").Value = 0.210715743495333
$ws.Range("J4").Value = 0.2107157434953329
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 114.5606336666667
$ws.Range("N4").Value = 343.681901
$ws.Range("O4").Value = 0.7368570786832789
$ws.Range("P4").Value = 0.736857078683279
$ws.Range("Q4").Value = 440.2389110302933
$ws.Range("R4").Value = 3962.150199272639
$ws.Range("S4").Value = 0.1552673871845462
$ws.Range("T4").Value = 0.1552673871845462

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl5"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.842846333333334
$ws.Range("H5").Value = 11.528539
$ws.Range("I5").Value = 0.210715743495333
$ws.Range("J5").Value = 0.2107157434953329
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 40.91135266666667
$ws.Range("N5").Value = 122.734058
$ws.Range("O5").Value = 0.263142921316721
$ws.Range("P5").Value = 0.2631429213167211
$ws.Range("Q5").Value = 157.2160415868069
$ws.Range("R5").Value = 1414.944374281262
$ws.Range("S5").Value = 0.05544835631078677
$ws.Range("T5").Value = 0.05544835631078677

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl5"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.228676
$ws.Range("H6").Value = 36.686028
$ws.Range("I6").Value = 0.6705380157807161
$ws.Range("J6").Value = 0.6705380157807161
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 114.5606336666667
$ws.Range("N6").Value = 343.681901
$ws.Range("O6").Value = 0.7368570786832789
$ws.Range("P6").Value = 0.736857078683279
$ws.Range("Q6").Value = 1400.924871464359
$ws.Range("R6").Value = 12608.32384317923
$ws.Range("S6").Value = 0.4940906834542608
$ws.Range("T6").Value = 0.4940906834542609

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl5"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.228676
$ws.Range("H7").Value = 36.686028
$ws.Range("I7").Value = 0.6705380157807161
$ws.Range("J7").Value = 0.6705380157807161
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 40.91135266666667
$ws.Range("N7").Value = 122.734058
$ws.Range("O7").Value = 0.263142921316721
$ws.Range("P7").Value = 0.2631429213167211
$ws.Range("Q7").Value = 500.2916764824026
$ws.Range("R7").Value = 4502.625088341624
$ws.Range("S7").Value = 0.1764473323264552
$ws.Range("T7").Value = 0.1764473323264553

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl5"
$ws.Range("C8").Value = "Ccr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4572023333333333
$ws.Range("H8").Value = 1.371607
$ws.Range("I8").Value = 0.02506988862928798
$ws.Range("J8").Value = 0.02506988862928798
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 114.5606336666667
$ws.Range("N8").Value = 343.681901
$ws.Range("O8").Value = 0.7368570786832789
$ws.Range("P8").Value = 0.736857078683279
$ws.Range("Q8").Value = 52.37738902054522
$ws.Range("R8").Value = 471.396501184907
$ws.Range("S8").Value = 0.0184729248982923
$ws.Range("T8").Value = 0.0184729248982923

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl5"
$ws.Range("C9").Value = "Ccr1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4572023333333333
$ws.Range("H9").Value = 1.371607
$ws.Range("I9").Value = 0.02506988862928798
$ws.Range("J9").Value = 0.02506988862928798
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 40.91135266666667
$ws.Range("N9").Value = 122.734058
$ws.Range("O9").Value = 0.263142921316721
$ws.Range("P9").Value = 0.2631429213167211
$ws.Range("Q9").Value = 18.70476589902288
$ws.Range("R9").Value = 168.342893091206
$ws.Range("S9").Value = 0.006596963730995687
$ws.Range("T9").Value = 0.006596963730995688
